# Refresh market-price-derived columns (H:N) on each profession sheet.
# Values below mirror a scheduled market-data recompute; target cells/rows
# were identified from the authoritative worksheet diff.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 23081072
$ws.Range("I40").Value = 5871.5
$ws.Range("J40").Value = 42859816
$ws.Range("K40").Value = 5871.5
$ws.Range("L40").Value = 42859816
$ws.Range("M40").Value = -5696.5
$ws.Range("N40").Value = -42860166
$ws.Range("H62").Value = 9499.777
$ws.Range("I62").Value = 8812.25
$ws.Range("K62").Value = 8812.25
$ws.Range("M62").Value = -8188.25
$ws.Range("H65").Value = 9499.777
$ws.Range("I65").Value = 8812.25
$ws.Range("K65").Value = 44061.25
$ws.Range("M65").Value = -40941.25
$ws.Range("H74").Value = 4731.8823
$ws.Range("I74").Value = 4620.1665
$ws.Range("K74").Value = 4620.1665
$ws.Range("M74").Value = -3684.1665
$ws.Range("H77").Value = 4731.8823
$ws.Range("I77").Value = 4620.1665
$ws.Range("K77").Value = 23100.8325
$ws.Range("M77").Value = -18420.8325
$ws.Range("H88").Value = 2394.6667
$ws.Range("J88").Value = 3489.5
$ws.Range("L88").Value = 3489.5
$ws.Range("N88").Value = -4301.5
$ws.Range("H91").Value = 2394.6667
$ws.Range("J91").Value = 3489.5
$ws.Range("L91").Value = 3489.5
$ws.Range("N91").Value = -6297.5
$ws.Range("H98").Value = 532757.4
$ws.Range("I98").Value = 3978.25
$ws.Range("J98").Value = 1439235.9
$ws.Range("K98").Value = 3978.25
$ws.Range("L98").Value = 1439235.9
$ws.Range("M98").Value = -2480.25
$ws.Range("N98").Value = -1442231.9
$ws.Range("H106").Value = 78174.62
$ws.Range("I106").Value = 84587.086
$ws.Range("K106").Value = 84587.086
$ws.Range("M106").Value = -83956.086
$ws.Range("H107").Value = 2388.439
$ws.Range("J107").Value = 3819
$ws.Range("L107").Value = 3819
$ws.Range("N107").Value = -7659
$ws.Range("H113").Value = 27916
$ws.Range("I113").Value = 28333.166
$ws.Range("J113").Value = 26664.5
$ws.Range("K113").Value = 28333.166
$ws.Range("L113").Value = 26664.5
$ws.Range("M113").Value = -25079.166
$ws.Range("N113").Value = -33172.5
$ws.Range("H122").Value = 532757.4
$ws.Range("I122").Value = 3978.25
$ws.Range("J122").Value = 1439235.9
$ws.Range("K122").Value = 11934.75
$ws.Range("L122").Value = 4317707.699999999
$ws.Range("M122").Value = -9484.75
$ws.Range("N122").Value = -4322607.699999999
$ws.Range("H137").Value = 79582.60000000001
$ws.Range("I137").Value = 95594.914
$ws.Range("J137").Value = 15533.333
$ws.Range("K137").Value = 286784.742
$ws.Range("L137").Value = 46599.999
$ws.Range("M137").Value = -284234.742
$ws.Range("N137").Value = -51699.999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2125.7646
$ws.Range("J102").Value = 3413.125
$ws.Range("L102").Value = 3413.125
$ws.Range("N102").Value = -6657.125
$ws.Range("H132").Value = 7985.381
$ws.Range("I132").Value = 4947.375
$ws.Range("K132").Value = 14842.125
$ws.Range("M132").Value = -12312.125

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2186.0789
$ws.Range("I134").Value = 1781.0741
$ws.Range("K134").Value = 5343.2223
$ws.Range("M134").Value = -2808.2223

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2605.8333
$ws.Range("I3").Value = 1410
$ws.Range("K3").Value = 1410
$ws.Range("M3").Value = -1297
$ws.Range("H16").Value = 2073
$ws.Range("I16").Value = 1942
$ws.Range("J16").Value = 2433.25
$ws.Range("K16").Value = 1942
$ws.Range("L16").Value = 2433.25
$ws.Range("M16").Value = -1655
$ws.Range("N16").Value = -3007.25
$ws.Range("H22").Value = 345.34482
$ws.Range("I22").Value = 319.54544
$ws.Range("J22").Value = 426.42856
$ws.Range("K22").Value = 319.54544
$ws.Range("L22").Value = 426.42856
$ws.Range("M22").Value = 30.45456000000001
$ws.Range("N22").Value = -1126.42856
$ws.Range("H31").Value = 828849.5600000001
$ws.Range("I31").Value = 1429988.4
$ws.Range("K31").Value = 1429988.4
$ws.Range("M31").Value = -1429693.4
$ws.Range("H34").Value = 828849.5600000001
$ws.Range("I34").Value = 1429988.4
$ws.Range("K34").Value = 1429988.4
$ws.Range("M34").Value = -1429786.4
$ws.Range("H62").Value = 4292.8237
$ws.Range("I62").Value = 4043.818
$ws.Range("K62").Value = 4043.818
$ws.Range("M62").Value = -3419.818
$ws.Range("H65").Value = 4292.8237
$ws.Range("I65").Value = 4043.818
$ws.Range("K65").Value = 20219.09
$ws.Range("M65").Value = -17099.09
$ws.Range("H81").Value = 44999
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 44999
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 2073
$ws.Range("I113").Value = 1942
$ws.Range("J113").Value = 2433.25
$ws.Range("K113").Value = 1942
$ws.Range("L113").Value = 2433.25
$ws.Range("M113").Value = 228
$ws.Range("N113").Value = -6773.25
$ws.Range("H122").Value = 2175.8823
$ws.Range("I122").Value = 1856.2307
$ws.Range("K122").Value = 5568.6921
$ws.Range("M122").Value = -3118.6921

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2324
$ws.Range("I39").Value = 998
$ws.Range("J39").Value = 2513.4285
$ws.Range("K39").Value = 2994
$ws.Range("L39").Value = 7540.2855
$ws.Range("M39").Value = -2700
$ws.Range("N39").Value = -8128.2855
$ws.Range("H40").Value = 124.75
$ws.Range("I40").Value = 116.333336
$ws.Range("J40").Value = 150
$ws.Range("K40").Value = 465.333344
$ws.Range("L40").Value = 600
$ws.Range("M40").Value = -396.333344
$ws.Range("N40").Value = -738
$ws.Range("H46").Value = 13750322
$ws.Range("I46").Value = 419.5
$ws.Range("K46").Value = 1258.5
$ws.Range("M46").Value = -1167.5
$ws.Range("H114").Value = 1474.4375
$ws.Range("I114").Value = 1055
$ws.Range("K114").Value = 3165
$ws.Range("M114").Value = 89
$ws.Range("H121").Value = 3045.923
$ws.Range("I121").Value = 1925.2
$ws.Range("J121").Value = 3746.375
$ws.Range("K121").Value = 5775.6
$ws.Range("L121").Value = 11239.125
$ws.Range("M121").Value = -4465.6
$ws.Range("N121").Value = -13859.125
$ws.Range("H129").Value = 3072.4443
$ws.Range("I129").Value = 3269.6667
$ws.Range("J129").Value = 2973.8333
$ws.Range("K129").Value = 9809.000100000001
$ws.Range("L129").Value = 8921.499899999999
$ws.Range("M129").Value = -4809.000100000001
$ws.Range("N129").Value = -18921.4999
$ws.Range("H130").Value = 2849.8
$ws.Range("I130").Value = 2849.8
$ws.Range("K130").Value = 8549.400000000001
$ws.Range("M130").Value = -3529.400000000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 253.59459
$ws.Range("I2").Value = 110.07692
$ws.Range("J2").Value = 331.33334
$ws.Range("K2").Value = 110.07692
$ws.Range("L2").Value = 331.33334
$ws.Range("M2").Value = 2.923079999999999
$ws.Range("N2").Value = -557.33334

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6462.884
$ws.Range("I7").Value = 6455.722
$ws.Range("K7").Value = 6455.722
$ws.Range("M7").Value = -6343.722
$ws.Range("H16").Value = 4271.7
$ws.Range("I16").Value = 576.3889
$ws.Range("J16").Value = 37529.5
$ws.Range("K16").Value = 576.3889
$ws.Range("L16").Value = 37529.5
$ws.Range("M16").Value = -406.3889
$ws.Range("N16").Value = -37869.5
$ws.Range("H40").Value = 9792.457
$ws.Range("I40").Value = 10842.759
$ws.Range("J40").Value = 4716
$ws.Range("K40").Value = 10842.759
$ws.Range("L40").Value = 4716
$ws.Range("M40").Value = -10706.759
$ws.Range("N40").Value = -4988
$ws.Range("H61").Value = 1250.4
$ws.Range("I61").Value = 1188.375
$ws.Range("K61").Value = 1188.375
$ws.Range("M61").Value = -986.375
$ws.Range("H68").Value = 2337.6978
$ws.Range("I68").Value = 2293.3513
$ws.Range("J68").Value = 2611.1667
$ws.Range("K68").Value = 2293.3513
$ws.Range("L68").Value = 2611.1667
$ws.Range("M68").Value = -1544.3513
$ws.Range("N68").Value = -4109.1667
$ws.Range("H71").Value = 2337.6978
$ws.Range("I71").Value = 2293.3513
$ws.Range("J71").Value = 2611.1667
$ws.Range("K71").Value = 11466.7565
$ws.Range("L71").Value = 13055.8335
$ws.Range("M71").Value = -7722.7565
$ws.Range("N71").Value = -20543.8335
$ws.Range("H100").Value = 2163.375
$ws.Range("I100").Value = 2199.6
$ws.Range("K100").Value = 2199.6
$ws.Range("M100").Value = -1658.6
$ws.Range("H113").Value = 1250.4
$ws.Range("I113").Value = 1188.375
$ws.Range("K113").Value = 1188.375
$ws.Range("M113").Value = 981.625
$ws.Range("H126").Value = 6462.884
$ws.Range("I126").Value = 6455.722
$ws.Range("K126").Value = 19367.166
$ws.Range("M126").Value = -16897.166

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 896.2
$ws.Range("I107").Value = 893.5
$ws.Range("K107").Value = 2680.5
$ws.Range("M107").Value = -760.5
$ws.Range("H132").Value = 3036.3447
$ws.Range("I132").Value = 4398.2856
$ws.Range("K132").Value = 13194.8568
$ws.Range("M132").Value = -10664.8568
$ws.Range("H136").Value = 165261.47
$ws.Range("I136").Value = 247650.64
$ws.Range("K136").Value = 742951.92
$ws.Range("M136").Value = -740401.92

Write-Output "Applied scheduled market-data refresh across all profession sheets."
